# Insert a new first column ("ID") with row-label identifiers, shifting
# the existing columns A:E (A, B, C, D, F headers) one column to the
# right so they become B:F.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data (columns A:E -> B:F) by inserting a blank column
# at the front.
$ws.Columns.Item(1).Insert()

# New header for the inserted column.
$ws.Range("A1").Value = "ID"

# Match the bold/centered header formatting used by the rest of row 1.
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)

# Row identifiers for the inserted ID column (rows 2-25).
$ids = @(
    "Hb 2", "Hb 3", "S 24", "S 28", "Hb 107", "Hb 66", "Hb 69", "Hb 95",
    "Hb 99", "Hb 92", "Hb 40", "Hb 41", "S 11", "Hb 57", "S 21", "S 22",
    "S 3", "S 4", "S 5", "Hb 74", "Hb 79", "Hb 32", "S 15", "S 16"
)

for ($i = 0; $i -lt $ids.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $ids[$i]
}
